$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C10: change the "Integer min" value for the R30 rule from 18 to 1
$ws.Range("C10").Value = 1
